$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 and row 5 swap their CODIGO (A), PRECIO (C) and Cant. Img (D) values.
# STATUS WEB (B) is "Disponible" in both rows already, so it is left as-is.
#
# Columns A and C hold text that looks numeric/currency ("10930745010",
# "$ 182.990"), so a leading apostrophe is used to force them to stay text
# cells instead of being auto-converted to numbers/currency by Excel.
$ws.Range("A3").Formula = '''10930745010'
$ws.Range("C3").Formula = '''$ 182.990'
$ws.Range("D3").Value2 = 11

$ws.Range("A5").Formula = '''10962389016'
$ws.Range("C5").Formula = '''$ 165.990'
$ws.Range("D5").Value2 = 3
